# Adds the new player "Yoe Mama" as row 4 of the roster, mirroring the
# layout/format already used by rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from row 3 into row 4 for all the columns we are
# --- about to fill in, so the new cells pick up the existing styles
# --- (centered hyperlink font, centered date format, etc.) instead of
# --- Excel inventing brand-new style records.
$ws.Range("B3:D3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F3:M3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Fill in the new player's data ---
$ws.Range("B4").Value = "Yoe Mama"
$ws.Range("C4").Value = "YoeMama@example.com"
$ws.Range("D4").Value = "GoPro123"
$ws.Range("E4").Value = "2/4/2000"
$ws.Range("F4").Value = "ABCD123456HDEFLL09"
$ws.Range("G4").Value = " 1111 Consit SS"
$ws.Range("H4").Value = "222-1234"
$ws.Range("I4").Value = "Yoe"
$ws.Range("J4").Value = "Mama"
$ws.Range("K4").Value = "Reyes"
$ws.Range("L4").Value = 1234567890
$ws.Range("M4").Value = 3

# --- Hyperlink the e-mail address, matching C2/C3 ---
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:YoeMama@example.com") | Out-Null
# Adding the hyperlink resets the cell's alignment away from the rest of
# the column, so put it back the way row 2/3 have it.
$ws.Range("C4").HorizontalAlignment = -4108   # xlCenter

# --- Match the final on-screen selection recorded in the sheet ---
$ws.Range("E15").Select() | Out-Null
